# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by a fresh handback run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 341f3963-...md
$overview.Range("G4").Value = "2016-08-23 18:49:11"

# zh-cn sheet (row for 341f3963-...md): Handoff / Handback datetimes
$zhcn.Range("H4").Value = "2016-08-23 18:48:58"
$zhcn.Range("K4").Value = "2016-08-23 18:49:32"

# de-de sheet (row for 341f3963-...md): Handoff datetime mirrors Overview's
# "Latest HO Xliff Generate Date", Handback datetime is its own value.
$dede.Range("H4").Value = "2016-08-23 18:49:11"
$dede.Range("K4").Value = "2016-08-23 18:49:43"
